$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "531×6=" "802×4="
Replace-Text "174×2=" "196×3="
Replace-Text "571×8=" "771×3="
Replace-Text "754×4=" "157×3="
Replace-Text "331×3=" "241×6="
Replace-Text "957×4=" "187×7="
Replace-Text "505×7=" "797×7="
Replace-Text "842×4=" "813×8="
Replace-Text "301×5=" "229×4="
Replace-Text "633×6=" "758×5="
Replace-Text "112×3=" "377×9="
Replace-Text "481×4=" "597×5="
Replace-Text "858×5=" "547×3="
Replace-Text "129×5=" "279×3="
Replace-Text "462×8=" "507×8="
Replace-Text "658×4=" "659×6="
Replace-Text "145×8=" "436×9="
Replace-Text "134×9=" "426×9="
Replace-Text "149×7=" "645×6="
Replace-Text "545×4=" "269×9="
Replace-Text "212×2=" "682×7="
Replace-Text "860×6=" "518×8="
Replace-Text "738×5=" "712×4="
Replace-Text "415×5=" "268×2="
Replace-Text "369×5=" "774×6="

Write-Output "Done"
